$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "37.572.95"
$ws.Range("E2").Value2 = "  +0.62%  "
$ws.Range("D3").Value2 = "2.017.99"
$ws.Range("E3").Value2 = "  +0.55%  "
$ws.Range("E4").Value2 = "  +0.02%  "
$ws.Range("D5").Value2 = "'262.77"
$ws.Range("E5").Value2 = "  +6.09%  "
$ws.Range("D6").Value2 = "'0.618"
$ws.Range("E6").Value2 = "  -2.05%  "
$ws.Range("D7").Value2 = "'1.00"
$ws.Range("E7").Value2 = "  +0.02%  "
$ws.Range("D8").Value2 = "'55.83"
$ws.Range("E8").Value2 = "  -7.65%  "
$ws.Range("D9").Value2 = "'0.384"
$ws.Range("E9").Value2 = "  +0.34%  "
$ws.Range("D10").Value2 = "'0.0777"
$ws.Range("E10").Value2 = "  -3.61%  "
$ws.Range("E11").Value2 = "  -1.93%  "
$ws.Range("B12").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value2 = "2.320.41"
$ws.Range("E12").Value2 = "  +0.89%  "
$ws.Range("B13").Value2 = "Chainlink"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value2 = "'14.35"
$ws.Range("E13").Value2 = "  -5.51%  "
$ws.Range("D14").Value2 = "'0.803"
$ws.Range("E14").Value2 = "  -5.42%  "
$ws.Range("D15").Value2 = "'20.68"
$ws.Range("E15").Value2 = "  -9.03%  "
$ws.Range("D16").Value2 = "'5.24"
$ws.Range("E16").Value2 = "  -4.08%  "
$ws.Range("D17").Value2 = "2.035.61"
$ws.Range("E17").Value2 = "  +1.43%  "
$ws.Range("D18").Value2 = "37.464.07"
$ws.Range("E18").Value2 = "  +0.54%  "
$ws.Range("D19").Value2 = "'69.69"
$ws.Range("E19").Value2 = "  -1.17%  "
$ws.Range("D20").Value2 = "0.0₃0840"
$ws.Range("E20").Value2 = "  -3.15%  "
$ws.Range("D21").Value2 = "'5.15"
$ws.Range("E21").Value2 = "  -1.09%  "
$ws.Range("D22").Value2 = "'228.01"
$ws.Range("E22").Value2 = "  -1.34%  "
$ws.Range("D23").Value2 = "'2.69"
$ws.Range("E23").Value2 = "  +7.32%  "
$ws.Range("E24").Value2 = "  -0.06%  "
$ws.Range("D25").Value2 = "'2.32"
$ws.Range("E25").Value2 = "  -1.81%  "
$ws.Range("D26").Value2 = "'163.56"
$ws.Range("E26").Value2 = "  -0.23%  "
$ws.Range("D27").Value2 = "'8.94"
$ws.Range("E27").Value2 = "  -5.18%  "
$ws.Range("D28").Value2 = "'19.66"
$ws.Range("E28").Value2 = "  -0.28%  "
$ws.Range("D29").Value2 = "'0.128"
$ws.Range("E29").Value2 = "  -12.25%  "
$ws.Range("E30").Value2 = "  +0.33%  "
$ws.Range("E31").Value2 = "  -1.55%  "
$ws.Range("D32").Value2 = "'0.0650"
$ws.Range("E32").Value2 = "  -0.37%  "
$ws.Range("D33").Value2 = "'4.62"
$ws.Range("E33").Value2 = "  -4.35%  "
$ws.Range("D34").Value2 = "'4.49"
$ws.Range("E34").Value2 = "  -1.17%  "
$ws.Range("D35").Value2 = "'2.38"
$ws.Range("E35").Value2 = "  +0.21%  "
$ws.Range("D36").Value2 = "'1.83"
$ws.Range("E36").Value2 = "  +1.07%  "
$ws.Range("E37").Value2 = "  +0.07%  "
$ws.Range("D38").Value2 = "'3.33"
$ws.Range("E38").Value2 = "  +1.45%  "
$ws.Range("D39").Value2 = "'5.21"
$ws.Range("E39").Value2 = "  -5.74%  "
$ws.Range("E40").Value2 = "  +4.63%  "
$ws.Range("D41").Value2 = "'1.21"
$ws.Range("E41").Value2 = "  +2.33%  "
$ws.Range("D42").Value2 = "'0.0941"
$ws.Range("E42").Value2 = "  -3.94%  "
$ws.Range("D43").Value2 = "'0.0213"
$ws.Range("E43").Value2 = "  -1.09%  "
$ws.Range("D44").Value2 = "1.402.44"
$ws.Range("E44").Value2 = "  +2.03%  "
$ws.Range("D45").Value2 = "'90.10"
$ws.Range("E45").Value2 = "  -0.94%  "
$ws.Range("D46").Value2 = "'15.64"
$ws.Range("E46").Value2 = "  -6.59%  "
$ws.Range("E47").Value2 = "  -2.03%  "
$ws.Range("D48").Value2 = "'7.07"
$ws.Range("E48").Value2 = "  -2.78%  "
$ws.Range("E49").Value2 = "  +0.80%  "
$ws.Range("B50").Value2 = "RocketPoolETH"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value2 = "2.211.01"
$ws.Range("E50").Value2 = "  +0.79%  "
$ws.Range("B51").Value2 = "NEARProtocol"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value2 = "'1.97"
$ws.Range("E51").Value2 = "  -2.13%  "
